$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (column H1) onto the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (row 1) - new columns I0 and IF
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-7 for columns I and J
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 8
